# Automatische test-sync: 2025-06-19 13:30:10
# Adds a new mail-log row (row 12) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows, and bumps the
# "Overig" tally on the "Dashboard" sheet from 5 to 6.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A12").Value = "Vragen over samenwerking"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D12").Value = "Overig"
$logs.Range("F12").Value = "2025-06-19 13:28:11"
$logs.Range("G12").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include row 12 ---
$catConditions = $logs.Range("D2:D11").FormatConditions
foreach ($i in 1..$catConditions.Count) {
    $catConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D12"))
}

$answeredConditions = $logs.Range("G2:G11").FormatConditions
foreach ($i in 1..$answeredConditions.Count) {
    $answeredConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G12"))
}

# --- Dashboard sheet: bump the "Overig" count from 5 to 6 ------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 6
